$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the new rows are treated as text so values like "003003" and
# "$10.75" are preserved verbatim instead of being coerced to numbers.
$ws.Range("A3:E8").NumberFormat = "@"

# Row 3
$ws.Range("A3").Value = "003003"
$ws.Range("B3").Value = "Red Jacket - Fuji Apple"
$ws.Range("C3").Value = "2"
$ws.Range("D3").Value = "$10.75"
$ws.Range("E3").Value = "$21.50"

# Row 4
$ws.Range("A4").Value = "003004"
$ws.Range("B4").Value = "Red Jacket - Rasp/Apple"
$ws.Range("C4").Value = "2"
$ws.Range("D4").Value = "$10.75"
$ws.Range("E4").Value = "$21.50"

# Row 5
$ws.Range("A5").Value = "003005"
$ws.Range("B5").Value = "Red Jacket - Strawberry (12oz)"
$ws.Range("C5").Value = "2"
$ws.Range("D5").Value = "$10.75"
$ws.Range("E5").Value = "$21.50"

# Row 6
$ws.Range("A6").Value = "003014"
$ws.Range("B6").Value = "Red Jacket - Fuji Apple 32oz"
$ws.Range("C6").Value = "1"
$ws.Range("D6").Value = "$24.85"
$ws.Range("E6").Value = "$24.85"

# Row 7
$ws.Range("A7").Value = "004020"
$ws.Range("B7").Value = "Natalie's - Orange Juice"
$ws.Range("C7").Value = "1"
$ws.Range("D7").Value = "$27.75"
$ws.Range("E7").Value = "$27.75"

# Row 8
$ws.Range("A8").Value = "004060"
$ws.Range("B8").Value = "Natalie's - Honey Tangarine"
$ws.Range("C8").Value = "1"
$ws.Range("D8").Value = "$14.35"
$ws.Range("E8").Value = "$14.35"
